# Rotate the weekly records among rows 2, 3, 4 and 9:
#   row2 <- old row9, row3 <- old row2, row4 <- old row3, row9 <- old row4
# Columns involved: D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "J", "K", "L", "M", "P")

# Capture the current ("before") values for the four affected rows.
$rows = @(2, 3, 4, 9)
$original = @{}
foreach ($r in $rows) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $vals
}

# Mapping: destination row -> source row (values to copy from the "before" snapshot)
$mapping = @{
    2 = 9
    3 = 2
    4 = 3
    9 = 4
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value = $original[$src][$c]
    }
}
